$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tidsredovisning")

# First move the "Summa" totals row from row 17 down to row 20,
# and update the SUM range to include the new rows (B2:B17).
$ws.Range("A20").Value = "Summa"
$ws.Range("B20").Formula = "=SUM(B2:B17)"

# Reuse date formatting (style) from the existing date column by
# copying formats from A16 before writing the new date serials.
$ws.Range("A16").Copy()
$ws.Range("A17:A19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 17: Spara ny post(tester)
$ws.Range("A17").Value = 45317
$ws.Range("B17").Value = 1
$ws.Range("C17").Value = "Spara ny post(tester)"

# Row 18: Kontrollera indata
$ws.Range("A18").Value = 45317
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = "Kontrollera indata"

# Row 19: Hämta enskild uppgift
$ws.Range("A19").Value = 45317
$ws.Range("B19").Value = 2
$ws.Range("C19").Value = "Hämta enskild uppgift"

# Update selection to match the new active cell
$ws.Range("C19").Select() | Out-Null
